# Daily attendance processing - 2025-10-30 19:19:08
# Normalises the "Recorded By" (column G) cell text so the contributor list
# always reads "System" first, followed by the human editor(s) - matching
# the canonical ordering produced by the attendance exporter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "Recorded By" values. Wherever a cell's comma-
# separated contributor list is exactly two names long and includes a
# human editor (dnasr281@gmail.com / admin@admin.com) alongside "System"
# (or the two human editors together), swap the order of the two names.
$dim = $ws.UsedRange
$lastRow = $dim.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $text = $cell.Value()

    if ($text -eq $null) { continue }

    $parts = $text -split ', '
    if ($parts.Count -ne 2) { continue }

    $first = $parts[0]
    $second = $parts[1]
    $humans = @('dnasr281@gmail.com', 'admin@admin.com')

    $firstIsHuman = $humans -contains $first
    $secondIsHuman = $humans -contains $second

    $swap = $false
    if ($firstIsHuman -and $second -eq 'System') { $swap = $true }
    if ($firstIsHuman -and $secondIsHuman) { $swap = $true }

    if ($swap) {
        $cell.Value = "$second, $first"
    }
}
